{"js": "// Applies the per-cell text replacements described by the diff:\n// the header date, plus each '<factor>\u00d7<digit>=' multiplication-table cell.\n// Every 'before' string is unique in the document, so a plain search +\n// insertText(\"Replace\") round-trip is safe for each pair.\nconst replacements = [\n  [\"2024-01-21 Sunday\", \"2024-01-22 Monday\"],\n  [\"145\u00d78=\", \"367\u00d78=\"],\n  [\"638\u00d79=\", \"854\u00d78=\"],\n  [\"800\u00d78=\", \"664\u00d78=\"],\n  [\"612\u00d79=\", \"734\u00d78=\"],\n  [\"267\u00d75=\", \"601\u00d75=\"],\n  [\"713\u00d78=\", \"818\u00d77=\"],\n  [\"375\u00d78=\", \"233\u00d73=\"],\n  [\"874\u00d72=\", \"546\u00d72=\"],\n  [\"188\u00d73=\", \"967\u00d76=\"],\n  [\"458\u00d76=\", \"735\u00d74=\"],\n  [\"129\u00d77=\", \"562\u00d72=\"],\n  [\"341\u00d76=\", \"142\u00d75=\"],\n  [\"257\u00d76=\", \"636\u00d76=\"],\n  [\"142\u00d76=\", \"667\u00d73=\"],\n  [\"566\u00d74=\", \"619\u00d75=\"],\n  [\"904\u00d74=\", \"534\u00d72=\"],\n  [\"744\u00d75=\", \"521\u00d79=\"],\n  [\"750\u00d77=\", \"886\u00d79=\"],\n  [\"450\u00d79=\", \"840\u00d77=\"],\n  [\"498\u00d77=\", \"118\u00d72=\"],\n  [\"723\u00d77=\", \"305\u00d78=\"],\n  [\"978\u00d79=\", \"886\u00d76=\"],\n  [\"417\u00d77=\", \"102\u00d75=\"],\n  [\"984\u00d75=\", \"341\u00d75=\"],\n  [\"412\u00d79=\", \"357\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the per-cell text replacements described by the diff:\n# the header date, plus each '<factor>x<digit>=' multiplication-table cell.\n# Every 'before' string is unique in the document, so Find/Replace with\n# wdReplaceAll (bounded by MatchCase/MatchWholeWord) safely retargets just\n# that one run per pair, leaving every other run untouched.\n$d = $word.ActiveDocument\n\n$findTexts = @(\n    '2024-01-21 Sunday',\n    '145\u00d78=',\n    '638\u00d79=',\n    '800\u00d78=',\n    '612\u00d79=',\n    '267\u00d75=',\n    '713\u00d78=',\n    '375\u00d78=',\n    '874\u00d72=',\n    '188\u00d73=',\n    '458\u00d76=',\n    '129\u00d77=',\n    '341\u00d76=',\n    '257\u00d76=',\n    '142\u00d76=',\n    '566\u00d74=',\n    '904\u00d74=',\n    '744\u00d75=',\n    '750\u00d77=',\n    '450\u00d79=',\n    '498\u00d77=',\n    '723\u00d77=',\n    '978\u00d79=',\n    '417\u00d77=',\n    '984\u00d75=',\n    '412\u00d79='\n)\n$replaceTexts = @(\n    '2024-01-22 Monday',\n    '367\u00d78=',\n    '854\u00d78=',\n    '664\u00d78=',\n    '734\u00d78=',\n    '601\u00d75=',\n    '818\u00d77=',\n    '233\u00d73=',\n    '546\u00d72=',\n    '967\u00d76=',\n    '735\u00d74=',\n    '562\u00d72=',\n    '142\u00d75=',\n    '636\u00d76=',\n    '667\u00d73=',\n    '619\u00d75=',\n    '534\u00d72=',\n    '521\u00d79=',\n    '886\u00d79=',\n    '840\u00d77=',\n    '118\u00d72=',\n    '305\u00d78=',\n    '886\u00d76=',\n    '102\u00d75=',\n    '341\u00d75=',\n    '357\u00d76='\n)\n\nfor ($i = 0; $i -lt $findTexts.Count; $i++) {\n    $findText = $findTexts[$i]\n    $replaceText = $replaceTexts[$i]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue=1, wdReplaceAll=2\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n"}
